$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "'" + '58.556.32'

$ws.Range("E2").Value = '  -1.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "'" + '2.487.08'

$ws.Range("E3").Value = '  -0.97%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'" + '526.58'

$ws.Range("E5").Value = '  -2.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "'" + '134.03'

$ws.Range("E6").Value = '  -3.23%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.92%  '

$ws.Range("E9").Value = '  -1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "'" + '0.156'

$ws.Range("E10").Value = '  -1.88%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "'" + '5.43'

$ws.Range("E11").Value = '  +0.94%  '

$ws.Range("E12").Value = '  -1.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'" + '2.927.76'

$ws.Range("E13").Value = '  -1.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "'" + '58.453.63'

$ws.Range("E14").Value = '  -1.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "'" + '22.49'

$ws.Range("E15").Value = '  -3.41%  '

$ws.Range("E16").Value = '  -2.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'" + '2.487.38'

$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "'" + '10.96'

$ws.Range("E18").Value = '  -1.51%  '

$ws.Range("E19").Value = '  -1.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "'" + '322.16'

$ws.Range("E20").Value = '  -1.30%  '

$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "'" + '5.84'

$ws.Range("E22").Value = '  -1.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "'" + '64.43'

$ws.Range("E23").Value = '  -1.56%  '

$ws.Range("E24").Value = '  -2.64%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'" + '0.999'

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("E26").Value = '  -3.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'" + '7.49'

$ws.Range("E27").Value = '  -2.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "'" + '0.0' + [char]0x2083 + '0756'

$ws.Range("E28").Value = '  -3.09%  '

$ws.Range("E29").Value = '  -4.60%  '

$ws.Range("E30").Value = '  -4.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "'" + '165.83'

$ws.Range("E31").Value = '  -1.94%  '

$ws.Range("E32").Value = '  -5.48%  '

$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'" + '0.999'

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "'" + '18.28'

$ws.Range("E35").Value = '  -1.48%  '

$ws.Range("E36").Value = '  -8.61%  '

$ws.Range("E37").Value = '  -3.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "'" + '1.51'

$ws.Range("E38").Value = '  -3.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'" + '0.799'

$ws.Range("E39").Value = '  -3.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'" + '3.54'

$ws.Range("E40").Value = '  -3.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "'" + '278.91'

$ws.Range("E41").Value = '  -2.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "'" + '4.99'

$ws.Range("E42").Value = '  -5.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "'" + '127.51'

$ws.Range("E44").Value = '  -3.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "'" + '0.0915'

$ws.Range("E45").Value = '  -2.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "'" + '0.0498'

$ws.Range("E46").Value = '  -2.76%  '

$ws.Range("E47").Value = '  -2.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "'" + '17.32'

$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "'" + '1.747.71'

$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "'" + '0.975'

$ws.Range("E50").Value = '  -1.47%  '

$ws.Range("E51").Value = '  -2.13%  '
